$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.953.53"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.757.98"
$ws.Range("E3").Value = "  -3.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.58"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3753"
$ws.Range("E7").Value = "  -4.85%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3351"
$ws.Range("E8").Value = "  -4.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.84"
$ws.Range("E9").Value = "  -5.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.121"
$ws.Range("E10").Value = "  -6.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07117"
$ws.Range("E11").Value = "  -6.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.19"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.161"
$ws.Range("E14").Value = "  -6.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.138"
$ws.Range("E15").Value = "  -1.47%  "
$ws.Range("D16").Value = "1.760.74"
$ws.Range("E16").Value = "  -3.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001047"
$ws.Range("E17").Value = "  -5.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06571"
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "79.96"
$ws.Range("E19").Value = "  -6.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.86"
$ws.Range("E21").Value = "  -6.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.246"
$ws.Range("E22").Value = "  -5.75%  "
$ws.Range("D23").Value = "27.911.26"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.62"
$ws.Range("E24").Value = "  -9.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.391"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.07"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.65"
$ws.Range("E27").Value = "  -9.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.309"
$ws.Range("E28").Value = "  -11.23%  "
$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.279"
$ws.Range("E29").Value = "  -16.07%  "
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").Value = "1.960.15"
$ws.Range("E30").Value = "  -3.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.06"
$ws.Range("E31").Value = "  -4.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.028"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.744"
$ws.Range("E33").Value = "  -8.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08730"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.14"
$ws.Range("E35").Value = "  -9.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02331"
$ws.Range("E36").Value = "  -4.99%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06174"
$ws.Range("E37").Value = "  -6.48%  "
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6511"
$ws.Range("E38").Value = "  -7.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.101"
$ws.Range("E39").Value = "  -8.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2099"
$ws.Range("E40").Value = "  -6.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.209"
$ws.Range("E41").Value = "  -5.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.451"
$ws.Range("E42").Value = "  -10.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.030"
$ws.Range("E43").Value = "  -6.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.57"
$ws.Range("E45").Value = "  -7.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.835"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5985"
$ws.Range("E47").Value = "  -8.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.36"
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.995"
$ws.Range("E49").Value = "  -8.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07193"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.168"
$ws.Range("E51").Value = "  +0.10%  "
